# Trade #165 closed at 2026-02-18 00:44:49 - unknown UNKNOWN +0.000%
#
# This script applies the following changes to live_trading_results.xlsx:
#  1. "Summary" sheet: Total Trades 192 -> 193, Win Rate % 43.23 -> 43.01
#  2. "Strategy Status" sheet: MarketMaking row - Trades 80 -> 81, Win Rate % 43.75 -> 43.21
#  3. "All Trades" sheet:
#       - row 194 (Trade #193, MarketMaking) closes: Exit Price, Status, Capital After,
#         Exit Reason, Duration updated
#       - two new open trades appended as rows 223/224 (Trade #222 HighProbConvergence,
#         Trade #223 MarketMaking)
#  4. "HighProbConvergence" sheet: new open trade appended as row 30 (Trade #222)
#  5. "MarketMaking" sheet:
#       - row 82 (Trade #193) closes, mirroring the "All Trades" update
#       - new open trade appended as row 104 (Trade #223)

$wb = $excel.ActiveWorkbook

# Helper: write a literal "YYYY-MM-DD" / date-looking string into a cell
# without letting Excel auto-convert it into a real date serial number -
# the source data always stores these as plain text.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 193
$summary.Range("B9").Value = 43.01

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking is row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 81
$status.Range("G6").Value = 43.21

# ---------------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out trade #193 (MarketMaking) recorded on row 194
$allTrades.Cells.Item(194, 7).Value = 0.01            # G194 Exit Price
$allTrades.Cells.Item(194, 8).Value = "CLOSED"        # H194 Status
$allTrades.Cells.Item(194, 11).Value = 99.45999999999999  # K194 Capital After
$allTrades.Cells.Item(194, 12).Value = "early_exit"   # L194 Exit Reason
$allTrades.Cells.Item(194, 13).Value = 0.17           # M194 Duration (min)

# New row 223: Trade #222, HighProbConvergence, OPEN
$allTrades.Cells.Item(223, 1).Value = 222
Set-TextValue $allTrades.Cells.Item(223, 2) "2026-02-18"
$allTrades.Cells.Item(223, 3).Value = "00:44:43"
$allTrades.Cells.Item(223, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(223, 5).Value = "UP"
$allTrades.Cells.Item(223, 6).Value = 0.01
$allTrades.Cells.Item(223, 7).Value = ""
$allTrades.Cells.Item(223, 8).Value = "OPEN"
$allTrades.Cells.Item(223, 9).Value = 0
$allTrades.Cells.Item(223, 10).Value = 0
$allTrades.Cells.Item(223, 11).Value = 100.3223499536821
$allTrades.Cells.Item(223, 12).Value = ""
$allTrades.Cells.Item(223, 13).Value = 0
$allTrades.Cells.Item(223, 14).Value = 0
$allTrades.Cells.Item(223, 15).Value = 0
$allTrades.Cells.Item(223, 16).Value = 0.95
$allTrades.Cells.Item(223, 17).Value = "Mean reversion UP: price 9.74% below mean (z=-2.00)"

# New row 224: Trade #223, MarketMaking, OPEN
$allTrades.Cells.Item(224, 1).Value = 223
Set-TextValue $allTrades.Cells.Item(224, 2) "2026-02-18"
$allTrades.Cells.Item(224, 3).Value = "00:44:43"
$allTrades.Cells.Item(224, 4).Value = "MarketMaking"
$allTrades.Cells.Item(224, 5).Value = "UP"
$allTrades.Cells.Item(224, 6).Value = 0.01
$allTrades.Cells.Item(224, 7).Value = ""
$allTrades.Cells.Item(224, 8).Value = "OPEN"
$allTrades.Cells.Item(224, 9).Value = 0
$allTrades.Cells.Item(224, 10).Value = 0
$allTrades.Cells.Item(224, 11).Value = 99.45858346467946
$allTrades.Cells.Item(224, 12).Value = ""
$allTrades.Cells.Item(224, 13).Value = 0
$allTrades.Cells.Item(224, 14).Value = 0
$allTrades.Cells.Item(224, 15).Value = 0
$allTrades.Cells.Item(224, 16).Value = 0.6
$allTrades.Cells.Item(224, 17).Value = "Normal spread capture: 225 bps"

# ---------------------------------------------------------------------------
# 4. HighProbConvergence sheet - append new open trade as row 30
#    Column order differs from "All Trades":
#    A Trade#, B Date, C Time, D Strategy, E Side, F Entry, G Exit, H Status,
#    I P&L%, J P&L$, K Capital After, L Entry Slippage, M Exit Slippage,
#    N Confidence, O Entry Reason, P Exit Reason, Q Duration
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(30, 1).Value = 222
Set-TextValue $hpc.Cells.Item(30, 2) "2026-02-18"
$hpc.Cells.Item(30, 3).Value = "00:44:43"
$hpc.Cells.Item(30, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(30, 5).Value = "UP"
$hpc.Cells.Item(30, 6).Value = 0.01
$hpc.Cells.Item(30, 7).Value = ""
$hpc.Cells.Item(30, 8).Value = "OPEN"
$hpc.Cells.Item(30, 9).Value = 0
$hpc.Cells.Item(30, 10).Value = 0
$hpc.Cells.Item(30, 11).Value = 100.3223499536821
$hpc.Cells.Item(30, 12).Value = 0
$hpc.Cells.Item(30, 13).Value = 0
$hpc.Cells.Item(30, 14).Value = 0.95
$hpc.Cells.Item(30, 15).Value = "Mean reversion UP: price 9.74% below mean (z=-2.00)"
$hpc.Cells.Item(30, 16).Value = ""
$hpc.Cells.Item(30, 17).Value = 0

# ---------------------------------------------------------------------------
# 5. MarketMaking sheet
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out trade #193 recorded on row 82 (same column order as HighProbConvergence)
$mm.Cells.Item(82, 7).Value = 0.01               # G82 Exit Price
$mm.Cells.Item(82, 8).Value = "CLOSED"           # H82 Status
$mm.Cells.Item(82, 11).Value = 99.45999999999999 # K82 Capital After
$mm.Cells.Item(82, 16).Value = "early_exit"      # P82 Exit Reason
$mm.Cells.Item(82, 17).Value = 0.17              # Q82 Duration (min)

# New row 104: Trade #223, OPEN
$mm.Cells.Item(104, 1).Value = 223
Set-TextValue $mm.Cells.Item(104, 2) "2026-02-18"
$mm.Cells.Item(104, 3).Value = "00:44:43"
$mm.Cells.Item(104, 4).Value = "MarketMaking"
$mm.Cells.Item(104, 5).Value = "UP"
$mm.Cells.Item(104, 6).Value = 0.01
$mm.Cells.Item(104, 7).Value = ""
$mm.Cells.Item(104, 8).Value = "OPEN"
$mm.Cells.Item(104, 9).Value = 0
$mm.Cells.Item(104, 10).Value = 0
$mm.Cells.Item(104, 11).Value = 99.45858346467946
$mm.Cells.Item(104, 12).Value = 0
$mm.Cells.Item(104, 13).Value = 0
$mm.Cells.Item(104, 14).Value = 0.6
$mm.Cells.Item(104, 15).Value = "Normal spread capture: 225 bps"
$mm.Cells.Item(104, 16).Value = ""
$mm.Cells.Item(104, 17).Value = 0
